$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -9
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -4
$ws.Range("F11").Value = -8
